$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last refreshed" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 07:22"

# Rows shifted because Hungria's updated data moved it above Kazajistan
# in the ranking. Row 64 now holds Hungria's fresh data, row 65 now
# holds the data Kazajistan previously had (unchanged values, new row).
$ws.Range("A64").Value = "Hungria"
$ws.Range("B64").Value = 2443
$ws.Range("C64").Value = 60
$ws.Range("D64").Value = 458
$ws.Range("E64").Value = 1723
$ws.Range("F64").Value = 61
$ws.Range("G64").Value = 12
$ws.Range("H64").Value = 262

$ws.Range("A65").Value = "Kazajistan"
$ws.Range("B65").Value = 2416
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 604
$ws.Range("E65").Value = 1787
$ws.Range("F65").Value = 29
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 25

# Same pattern: Tanzania's updated data moved it above Jamaica.
$ws.Range("A123").Value = "Tanzania"
$ws.Range("B123").Value = 299
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 48
$ws.Range("E123").Value = 241
$ws.Range("F123").Value = 7
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 10

$ws.Range("A124").Value = "Jamaica"
$ws.Range("B124").Value = 288
$ws.Range("C124").Value = 31
$ws.Range("D124").Value = 28
$ws.Range("E124").Value = 253
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 7
